$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update "Semestre ideal" value from EA-7 to EA-9 (row 9)
$ws.Range("B9").Value = "EA-9"
$ws.Range("C9").Value = "EA-9"

# 2. Copy the formatting (style + row height) of the existing "Requisitos" block
#    (rows 24:26) down into the new rows (27:51) that will be added.
$ws.Range("B24:C26").Copy($ws.Range("B27:C29"))
$ws.Range("B24:C26").Copy($ws.Range("B30:C32"))
$ws.Range("B24:C26").Copy($ws.Range("B33:C35"))
$ws.Range("B24:C26").Copy($ws.Range("B36:C38"))
$ws.Range("B24:C26").Copy($ws.Range("B39:C41"))
$ws.Range("B24:C26").Copy($ws.Range("B42:C44"))
$ws.Range("B24:C26").Copy($ws.Range("B45:C47"))
$ws.Range("B24:C26").Copy($ws.Range("B48:C50"))
$ws.Range("B24:C24").Copy($ws.Range("B51:C51"))
$ws.Rows("27:51").RowHeight = 30

# 3. Replace the text of the 3 pre-existing requirement rows (now "Requisito fraco")
$ws.Range("B24").Value = "LOB1003 -  Cálculo I  (Requisito fraco)`n"
$ws.Range("C24").Value = "LOB1003 -  Cálculo I  (Requisito fraco)`n"
$ws.Range("B25").Value = "LOB1004 -  Cálculo II  (Requisito fraco)`n"
$ws.Range("C25").Value = "LOB1004 -  Cálculo II  (Requisito fraco)`n"
$ws.Range("B26").Value = "LOB1006 -  Cálculo IV  (Requisito fraco)`n"
$ws.Range("C26").Value = "LOB1006 -  Cálculo IV  (Requisito fraco)`n"

# 4. Fill the new rows 27-51 with the remaining "Requisito fraco" values
$ws.Range("B27").Value = "LOB1011 -  Eletricidade Aplicada  (Requisito fraco)`n"
$ws.Range("C27").Value = "LOB1011 -  Eletricidade Aplicada  (Requisito fraco)`n"
$ws.Range("B28").Value = "LOB1012 -  Estatística  (Requisito fraco)`n"
$ws.Range("C28").Value = "LOB1012 -  Estatística  (Requisito fraco)`n"
$ws.Range("B29").Value = "LOB1018 -  Física I  (Requisito fraco)`n"
$ws.Range("C29").Value = "LOB1018 -  Física I  (Requisito fraco)`n"
$ws.Range("B30").Value = "LOB1019 -  Física II  (Requisito fraco)`n"
$ws.Range("C30").Value = "LOB1019 -  Física II  (Requisito fraco)`n"
$ws.Range("B31").Value = "LOB1021 -  Física IV  (Requisito fraco)`n"
$ws.Range("C31").Value = "LOB1021 -  Física IV  (Requisito fraco)`n"
$ws.Range("B32").Value = "LOB1024 -  Mecânica  (Requisito fraco)`n"
$ws.Range("C32").Value = "LOB1024 -  Mecânica  (Requisito fraco)`n"
$ws.Range("B33").Value = "LOB1036 -  Geometria Analítica  (Requisito fraco)`n"
$ws.Range("C33").Value = "LOB1036 -  Geometria Analítica  (Requisito fraco)`n"
$ws.Range("B34").Value = "LOB1037 -  Àlgebra Linear  (Requisito fraco)`n"
$ws.Range("C34").Value = "LOB1037 -  Àlgebra Linear  (Requisito fraco)`n"
$ws.Range("B35").Value = "LOB1038 -  Física Experimental I  (Requisito fraco)`n"
$ws.Range("C35").Value = "LOB1038 -  Física Experimental I  (Requisito fraco)`n"
$ws.Range("B36").Value = "LOB1039 -  Física Experimental III  (Requisito fraco)`n"
$ws.Range("C36").Value = "LOB1039 -  Física Experimental III  (Requisito fraco)`n"
$ws.Range("B37").Value = "LOB1041 -  Física Experimental II  (Requisito fraco)`n"
$ws.Range("C37").Value = "LOB1041 -  Física Experimental II  (Requisito fraco)`n"
$ws.Range("B38").Value = "LOB1042 -  Física Experimental IV  (Requisito fraco)`n"
$ws.Range("C38").Value = "LOB1042 -  Física Experimental IV  (Requisito fraco)`n"
$ws.Range("B39").Value = "LOB1045 -  Leitura e Produção de Textos Acadêmicos  (Requisito fraco)`n"
$ws.Range("C39").Value = "LOB1045 -  Leitura e Produção de Textos Acadêmicos  (Requisito fraco)`n"
$ws.Range("B40").Value = "LOB1052 -  Cálculo III  (Requisito fraco)`n"
$ws.Range("C40").Value = "LOB1052 -  Cálculo III  (Requisito fraco)`n"
$ws.Range("B41").Value = "LOB1053 -  Física III  (Requisito fraco)`n"
$ws.Range("C41").Value = "LOB1053 -  Física III  (Requisito fraco)`n"
$ws.Range("B42").Value = "LOB1056 -  Introdução aos Métodos Numéricos e Computacionais  (Requisito fraco)`n"
$ws.Range("C42").Value = "LOB1056 -  Introdução aos Métodos Numéricos e Computacionais  (Requisito fraco)`n"
$ws.Range("B43").Value = "LOB1232 -  Licenciamento Ambiental  (Requisito fraco)`n"
$ws.Range("C43").Value = "LOB1232 -  Licenciamento Ambiental  (Requisito fraco)`n"
$ws.Range("B44").Value = "LOB1257 -  Sistema de Abastecimento e Tratamento de Água  (Requisito fraco)`n"
$ws.Range("C44").Value = "LOB1257 -  Sistema de Abastecimento e Tratamento de Água  (Requisito fraco)`n"
$ws.Range("B45").Value = "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito fraco)`n"
$ws.Range("C45").Value = "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito fraco)`n"
$ws.Range("B46").Value = "LOM3081 -  Introdução à Mecânica dos Sólidos  (Requisito fraco)`n"
$ws.Range("C46").Value = "LOM3081 -  Introdução à Mecânica dos Sólidos  (Requisito fraco)`n"
$ws.Range("B47").Value = "LOQ4095 -  Química Geral Experimental  (Requisito fraco)`n"
$ws.Range("C47").Value = "LOQ4095 -  Química Geral Experimental  (Requisito fraco)`n"
$ws.Range("B48").Value = "LOQ4097 -  Fundamentos de Química para Engenharia I (Requisito fraco)`n"
$ws.Range("C48").Value = "LOQ4097 -  Fundamentos de Química para Engenharia I (Requisito fraco)`n"
$ws.Range("B49").Value = "LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito fraco)`n"
$ws.Range("C49").Value = "LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito fraco)`n"
$ws.Range("B50").Value = "LOQ4233 -  Gestão de Negócios  (Requisito fraco)`n"
$ws.Range("C50").Value = "LOQ4233 -  Gestão de Negócios  (Requisito fraco)`n"
$ws.Range("B51").Value = "LOQ4247 -  Desenho Assistido por Computador  (Requisito fraco)`n"
$ws.Range("C51").Value = "LOQ4247 -  Desenho Assistido por Computador  (Requisito fraco)`n"

Write-Output "edit complete"
